$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# New "Hs" (heated space fraction) values for rows 2-19 on the ARCHITECTURE sheet
$values = @(0.25, 0.25, 0.84, 0.84, 0.84, 0.84, 0.84, 0.7, 0.67, 0.84, 0.67, 0, 1, 0, 1, 0.67, 0.67, 0.67)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

$ws.Range("E12").Select()
